# Weekly data refresh: a new week's price record is inserted at row 14,
# pushing all subsequent records (previously rows 14-59) down by one row
# (to rows 15-60). The sheet's used range grows from A1:R59 to A1:R60.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 14, shifting rows 14:59
# down to 15:60.
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with this week's record.
$ws.Range("A14").Value2 = 11
$ws.Range("B14").Value2 = "Vega Monumental Concepción"
$ws.Range("C14").Value2 = "Bíobío"
$ws.Range("D14").Value2 = 45133
$ws.Range("E14").Value2 = 8
$ws.Range("F14").Value2 = 100112026
$ws.Range("G14").Value2 = "Haba"
$ws.Range("H14").Value2 = "Sin especificar"
$ws.Range("I14").Value2 = "Primera"
$ws.Range("J14").Value2 = 100
$ws.Range("K14").Value2 = 15000
$ws.Range("L14").Value2 = 16000
$ws.Range("M14").Value2 = 15500
$ws.Range("N14").Value2 = "`$/saco 25 kilos"
$ws.Range("O14").Value2 = "Provincia de Limarí"
$ws.Range("P14").Value2 = 620
$ws.Range("Q14").Value2 = 25
$ws.Range("R14").Value2 = "Hortaliza"
